$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit cyclically permutes the data rows 3-6 (row 1 = header, row 2 is
# untouched): the record that used to be on row 5 moves to row 3, the one on
# row 3 moves to row 4, the one on row 6 moves to row 5, and the one on row 4
# moves to row 6. Concretely (matching the target OOXML diff) this means
# only the following per-cell values change; every other cell in rows 3-6
# (C, K, N, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AF, AG, AT, AW, AX, AY) is
# identical before/after the permutation and is left untouched.

$ws.Range("A3").Value = 111701910
$ws.Range("P3").Value = "Kyrkogården (Kyrkogården), Nrk"
$ws.Range("Q3").Value = 516978.9846792166
$ws.Range("R3").Value = 6574635.767148005

$ws.Range("A4").Value = 111701829
$ws.Range("B4").Value = 90687
$ws.Range("E4").Value = 5964
$ws.Range("F4").Value = "Fjällig taggsvamp s.str."
$ws.Range("G4").Value = "Sarcodon imbricatus s.str."
$ws.Range("H4").Value = "(L.:Fr.) P.Karst."
$ws.Range("P4").Value = "Myrövägen öster (Myrövägen öster), Nrk"
$ws.Range("Q4").Value = 516894.5773385105
$ws.Range("R4").Value = 6574639.474785783

$ws.Range("A5").Value = 111702271
$ws.Range("B5").Value = 90709
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 5448
$ws.Range("F5").Value = "Svartvit taggsvamp"
$ws.Range("G5").Value = "Phellodon connatus"
$ws.Range("H5").Value = "(Schultz) nom.prov"
# I5/J5 hold a free-text count ("3 fruktkroppar"), not numbers - prefix the
# numeric-looking "3" with an apostrophe so Excel stores it as text instead
# of silently converting it to the number 3.
$ws.Range("I5").Value = "'3"
$ws.Range("J5").Value = "fruktkroppar"
$ws.Range("Q5").Value = 516923.6092008445
$ws.Range("R5").Value = 6574666.663922376

$ws.Range("A6").Value = 111702281
$ws.Range("B6").Value = 89183
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 3215
$ws.Range("F6").Value = "Rödgul trumpetsvamp"
$ws.Range("G6").Value = "Craterellus lutescens"
$ws.Range("H6").Value = "(Fr.) Fr."
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("Q6").Value = 516917.5151204841
$ws.Range("R6").Value = 6574656.936104885
